$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 23 content
$ws.Range("A23").Value = 5
$ws.Range("B23").Value = $null
$ws.Range("C23").Value = "NPC02SXON-RC"
$ws.Range("D23").Value = "2 (1 x 2) Position Shunt Connector  Open Top, Grip 0.100`" (2.54mm) Gold"
$ws.Range("E23").Value = "NPC02SXON-RC"
$ws.Range("F23").Value = "Sullins Connector Solutions"

# Set row height for row 11 and row 23
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30

# Set selection and scroll
$ws.Range("B23").Select()
